$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 686, shifting existing rows 686-733 down to 688-735
$ws.Range("A686:A687").EntireRow.Insert()

# Populate new row 686 (Primera)
$ws.Cells.Item(686,1).Value = 8
$ws.Cells.Item(686,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(686,3).Value = "Coquimbo"
$ws.Cells.Item(686,4).Value = 44746
$ws.Cells.Item(686,5).Value = 4
$ws.Cells.Item(686,6).Value = 100112023
$ws.Cells.Item(686,7).Value = "Brócoli"
$ws.Cells.Item(686,8).Value = "Sin especificar"
$ws.Cells.Item(686,9).Value = "Primera"
$ws.Cells.Item(686,10).Value = 2540
$ws.Cells.Item(686,11).Value = 750
$ws.Cells.Item(686,12).Value = 800
$ws.Cells.Item(686,13).Value = 775
$ws.Cells.Item(686,14).Value = "`$/unidad"
$ws.Cells.Item(686,15).Value = "Provincia del Elquí"
$ws.Cells.Item(686,16).Value = 775
$ws.Cells.Item(686,17).Value = 1
$ws.Cells.Item(686,18).Value = "Hortaliza"

# Populate new row 687 (Segunda)
$ws.Cells.Item(687,1).Value = 8
$ws.Cells.Item(687,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(687,3).Value = "Coquimbo"
$ws.Cells.Item(687,4).Value = 44746
$ws.Cells.Item(687,5).Value = 4
$ws.Cells.Item(687,6).Value = 100112023
$ws.Cells.Item(687,7).Value = "Brócoli"
$ws.Cells.Item(687,8).Value = "Sin especificar"
$ws.Cells.Item(687,9).Value = "Segunda"
$ws.Cells.Item(687,10).Value = 1400
$ws.Cells.Item(687,11).Value = 650
$ws.Cells.Item(687,12).Value = 700
$ws.Cells.Item(687,13).Value = 675
$ws.Cells.Item(687,14).Value = "`$/unidad"
$ws.Cells.Item(687,15).Value = "Provincia del Elquí"
$ws.Cells.Item(687,16).Value = 675
$ws.Cells.Item(687,17).Value = 1
$ws.Cells.Item(687,18).Value = "Hortaliza"
